$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Auto date fields: 19-06-2021 -> 01-07-2021 (notes master) and
#    6/19/2021 -> 7/1/2021 (four slide layouts: Section Header,
#    Content with Caption, Picture with Caption, Vertical Title and Text).
# ---------------------------------------------------------------------------

function Update-DateField {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -eq "19-06-2021") {
                $shp.TextFrame.TextRange.Text = "01-07-2021"
            } elseif ($t -eq "6/19/2021") {
                $shp.TextFrame.TextRange.Text = "7/1/2021"
            }
        }
    }
}

# Notes master date placeholder
Update-DateField $p.NotesMaster.Shapes

# Slide-layout date placeholders (layouts are indexed off the slide master)
$sm = $p.SlideMaster
for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    Update-DateField $cl.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 3, shape "CustomShape 6": "Getting Started" -> "getting started"
#    (re-typed initial "G"/"S" as lower-case "g"/"s", keeping the rest of the
#    run's formatting). The shape uses spAutoFit, so its height is
#    recomputed automatically once the text is edited.
# ---------------------------------------------------------------------------

$slide3 = $p.Slides.Item(3)
$titleShape = $slide3.Shapes.Item(6)
$tr = $titleShape.TextFrame.TextRange

$tr.Characters(1, 1).Text = "g"
$tr.Characters(9, 1).Text = "s"
